# Refresh the "cryptos" price/volume snapshot (Price = column D, Volume(1h) = column E)
# for rows 2-51, mirroring the upstream coinranking.com scrape that feeds this sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading-apostrophe text prefix, so Excel keeps storing numeric-looking prices
# (e.g. "214.72") as literal Text - exactly how they were stored before the edit -
# instead of silently coercing them to the Number type.
$apostrophe = "'"

$ws.Range("D2").Value = '25.898.14'
$ws.Range("E2").Value = '  -0.55%  '
$ws.Range("D3").Value = '1.639.99'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("D4").Value = $apostrophe + '1.005'
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").Value = $apostrophe + '214.72'
$ws.Range("E5").Value = '  -0.23%  '
$ws.Range("D6").Value = $apostrophe + '0.5045'
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("E7").Value = '  -0.57%  '
$ws.Range("D8").Value = $apostrophe + '0.2564'
$ws.Range("E8").Value = '  -0.34%  '
$ws.Range("D9").Value = $apostrophe + '0.06380'
$ws.Range("E9").Value = '  -0.79%  '
$ws.Range("D10").Value = $apostrophe + '19.56'
$ws.Range("E10").Value = '  +0.49%  '
$ws.Range("D11").Value = $apostrophe + '0.07793'
$ws.Range("E11").Value = '  +0.94%  '
$ws.Range("D12").Value = '1.667.81'
$ws.Range("E12").Value = '  +1.61%  '
$ws.Range("D13").Value = $apostrophe + '4.271'
$ws.Range("E13").Value = '  +0.59%  '
$ws.Range("D14").Value = $apostrophe + '0.5418'
$ws.Range("E14").Value = '  -0.49%  '
$ws.Range("D15").Value = '0.0₅7861'
$ws.Range("E15").Value = '  -0.40%  '
$ws.Range("D16").Value = $apostrophe + '64.61'
$ws.Range("E16").Value = '  +1.56%  '
$ws.Range("D17").Value = '25.916.75'
$ws.Range("E17").Value = '  -0.43%  '
$ws.Range("E18").Value = '  -0.60%  '
$ws.Range("D19").Value = $apostrophe + '197.87'
$ws.Range("E19").Value = '  -3.00%  '
$ws.Range("D20").Value = $apostrophe + '4.380'
$ws.Range("E20").Value = '  +2.20%  '
$ws.Range("D21").Value = $apostrophe + '9.947'
$ws.Range("E21").Value = '  -0.34%  '
$ws.Range("D22").Value = $apostrophe + '5.985'
$ws.Range("E22").Value = '  +0.51%  '
$ws.Range("D23").Value = $apostrophe + '1.004'
$ws.Range("E23").Value = '  -0.57%  '
$ws.Range("D24").Value = $apostrophe + '1.861'
$ws.Range("E24").Value = '  -3.57%  '
$ws.Range("D25").Value = $apostrophe + '140.06'
$ws.Range("E25").Value = '  -0.88%  '
$ws.Range("D26").Value = $apostrophe + '0.1141'
$ws.Range("E26").Value = '  -1.08%  '
$ws.Range("D27").Value = $apostrophe + '6.832'
$ws.Range("E27").Value = '  +1.51%  '
$ws.Range("E28").Value = '  -0.12%  '
$ws.Range("D29").Value = $apostrophe + '1.240'
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").Value = $apostrophe + '0.04929'
$ws.Range("E30").Value = '  -2.42%  '
$ws.Range("D31").Value = $apostrophe + '3.257'
$ws.Range("E31").Value = '  +0.37%  '
$ws.Range("D32").Value = $apostrophe + '3.186'
$ws.Range("E32").Value = '  -0.15%  '
$ws.Range("D33").Value = $apostrophe + '1.527'
$ws.Range("E33").Value = '  -0.89%  '
$ws.Range("D34").Value = $apostrophe + '2.363'
$ws.Range("E34").Value = '  +0.96%  '
$ws.Range("D35").Value = $apostrophe + '0.8914'
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").Value = $apostrophe + '2.605'
$ws.Range("E36").Value = '  -0.88%  '
$ws.Range("D37").Value = '1.137.14'
$ws.Range("E37").Value = '  -0.65%  '
$ws.Range("D38").Value = $apostrophe + '0.5533'
$ws.Range("E38").Value = '  -1.50%  '
$ws.Range("D39").Value = $apostrophe + '0.01557'
$ws.Range("E39").Value = '  -0.95%  '
$ws.Range("D40").Value = $apostrophe + '1.002'
$ws.Range("E40").Value = '  -0.67%  '
$ws.Range("D41").Value = $apostrophe + '5.655'
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").Value = $apostrophe + '0.8179'
$ws.Range("E42").Value = '  +1.12%  '
$ws.Range("D43").Value = $apostrophe + '99.15'
$ws.Range("E43").Value = '  -0.71%  '
$ws.Range("D44").Value = '0.0₈124'
$ws.Range("E44").Value = '  +9.30%  '
$ws.Range("D45").Value = '1.778.30'
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").Value = $apostrophe + '0.4524'
$ws.Range("E46").Value = '  -0.14%  '
$ws.Range("D47").Value = $apostrophe + '55.15'
$ws.Range("E47").Value = '  +0.52%  '
$ws.Range("E48").Value = '  -0.43%  '
$ws.Range("E49").Value = '  +0.31%  '
$ws.Range("E50").Value = '  -0.21%  '
$ws.Range("D51").Value = $apostrophe + '0.09517'
$ws.Range("E51").Value = '  +2.21%  '
